{"js": "// Insert two new paragraphs (Client Name / Company Name) right after the\n// document title (Heading1) and before the \"1. \u0428\u0430\u0440\u0442\u043d\u043e\u043c\u0430 \u043f\u0440\u0435\u0434\u043c\u0435\u0442\u0438\" heading.\nconst body = context.document.body;\nconst allParas = body.paragraphs;\nallParas.load(\"items/text\");\nawait context.sync();\n\nconst titlePara = allParas.items[0];\nconst clientPara = titlePara.insertParagraph(\"Client Name: Whitney Hurley\", Word.InsertLocation.after);\nclientPara.styleBuiltIn = Word.Style.normal;\nconst companyPara = clientPara.insertParagraph(\"Company Name: Roberts and Richards Co\", Word.InsertLocation.after);\ncompanyPara.styleBuiltIn = Word.Style.normal;\nawait context.sync();\n\n// Re-load the paragraphs (indices shifted by the two inserts above) and\n// replace the body text of the \"2.2.\" and \"2.3.\" clauses in place so the\n// paragraphs themselves (and any other formatting) are preserved.\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet p22 = null;\nlet p23 = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (p22 === null && t.indexOf(\"2.2. \u0418\u043d\u0432\u0435\u0441\u0442\u043e\u0440 \u0439\u0438\u0433\u0438\u043c \u043c\u0438\u043a\u0434\u043e\u0440\u0438\u043d\u0438 \u0442\u0443\u043b\u0438\u043a \u0445\u0430\u0436\u043c\u0434\u0430\") === 0) {\n    p22 = paras.items[i];\n  } else if (p23 === null && t.indexOf(\"2.3. \u041c\u0430\u0437\u043a\u0443\u0440 \u0448\u0430\u0440\u0442\u043d\u043e\u043c\u0430\u0434\u0430 \u043a\u0423\u0440\u0441\u0430\u0442\u0438\u043b\u0433\u0430\u043d\") === 0) {\n    p23 = paras.items[i];\n  }\n}\n\nif (!p22 || !p23) {\n  throw new Error(\"Could not locate target paragraphs 2.2/2.3 to update\");\n}\n\np22.getRange().insertText(\"2.2. \u0418\u043d\u0432\u0435\u0441\u0442\u043e\u0440 \u0439\u0438\u0433\u0438\u043c \u043c\u0438\u043a\u0434\u043e\u0440\u0438\u043d\u0438 \u0442\u0443\u043b\u0438\u043a \u0445\u0430\u0436\u043c\u0434\u0430 \u0416\u0430\u043c\u0433\u0430\u0440\u043c\u0430 \u0445\u0438\u0441\u043e\u0431 \u0440\u0430\u043a\u0430\u043c\u0438\u0433\u0430 \u0443\u0442\u043a\u0430\u0437\u0438\u0448 \u0439\u0423\u043b\u0438 \u0431\u0438\u043b\u0430\u043d \u043a\u0443\u0439\u0438\u0434\u0430\u0433\u0438 \u0442\u0430\u0440\u0442\u0438\u0431\u0434\u0430 \u0442\u0423\u043b\u0430\u043d\u0430\u0434\u0438: \u2014 20 (\u0439\u0438\u0433\u0438\u0440\u043c\u0430) \u0444\u043e\u0438\u0437\u0438 432.480.000 (\u0442\u0443\u0440\u0442 \u044e\u0437 \u0443\u0442\u0442\u0438\u0437 \u0438\u043a\u043a\u0438 \u043c\u0438\u043b\u043b\u0438\u043e\u043d \u0442\u0443\u0440\u0442 \u044e\u0437 \u0441\u0430\u043a\u0441\u043e\u043d \u043c\u0438\u043d\u0433) \u0441\u0443\u043c \u043c\u0438\u043a\u0434\u043e\u0440\u0438\u0434\u0430\u0433\u0438 \u043e\u043b\u0434\u0438\u043d\u0434\u0430\u043d \u0442\u0443\u043b\u043e\u0432\u043d\u0438 \u043c\u0430\u0437\u043a\u0443\u0440 \u0448\u0430\u0440\u0442\u043d\u043e\u043c\u0430 \u0438\u043c\u0437\u043e\u043b\u0430\u043d\u0433\u0430\u043d \u0441\u0430\u043d\u0430\u0434\u0430\u043d \u0431\u043e\u0448\u043b\u0430\u0431 \u0417 (\u0443\u0447) \u0438\u0448 \u043a\u0443\u043d\u0438\u0434\u0430\u043d \u043a\u0435\u0447\u0438\u043a\u0442\u0438\u0440\u043c\u0430\u0433\u0430\u043d \u0445\u043e\u043b\u0434\u0430 \u0430\u043c\u0430\u043b\u0433\u0430 \u043e\u0448\u0438\u0440\u0430\u0434\u0438; \u2014 \u0439\u0438\u0433\u0438\u043c\u043d\u0438\u043d\u0433 \u043a\u043e\u043b\u0433\u0430\u043d 80 (\u0441\u0430\u043a\u0441\u043e\u043d) \u0444\u043e\u0438\u0437\u0438 1.729.920.000 (\u0431\u0438\u0440 \u043c\u0438\u043b\u043b\u0438\u0430\u0440\u0434 \u0439\u0435\u0442\u0442\u0438 \u044e\u0437 \u0439\u0438\u0433\u0438\u0440\u043c\u0430 \u0442\u0443\u043a\u043a\u0438\u0437 \u043c\u0438\u043b\u043b\u0438\u043e\u043d \u0442\u0443\u043a\u043a\u0438\u0437 \u044e\u0437 \u0439\u0438\u0433\u0438\u0440\u043c\u0430 \u043c\u0438\u043d\u0433) \u0441\u0443\u043c \u043c\u0438\u043a\u0434\u043e\u0440\u0434\u0430\u0433\u0438 \u043a\u044a\u0448\u0438\u043c\u0447\u0430 \u0442\u0443\u043b\u043e\u0432\u043d\u0438 \u043a\u0443\u043f \u043a\u0430\u0432\u0430\u0442\u043b\u0438 \u0442\u0443\u0440\u0430\u0440 \u0436\u043e\u0439 \u0431\u0438\u043d\u043e\u0441\u0438 (\u0438\u043d\u0448\u043e\u043e\u0442\u0438) \u043a\u0443\u0440\u0438\u0431 \u0431\u0438\u0442\u043a\u0430\u0437\u0438\u043b\u0433\u0443\u043d\u0433\u0430 \u043a\u0430\u0434\u0430\u0440, \u0430\u043c\u043c\u043e 2023 \u0439\u0438\u043b 1 \u043d\u043e\u044f\u0431\u0440\u0434\u0430\u043d \u043a\u0435\u0447\u0438\u043a\u0442\u0438\u0440\u043c\u0430\u0441\u0434\u0430\u043d \u043a\u0443\u0439\u0438\u0434\u0430\u0433\u0438 \u043c\u0443\u0445\u043b\u0430\u0442\u043b\u0430\u0440\u0434\u0430 \u0430\u043c\u0430\u043b\u0433\u0430 \u043e\u0448\u0438\u0440\u0430\u0434\u0438.\", Word.InsertLocation.replace);\np23.getRange().insertText(\"2.3. \u0418\u0438\u0433\u0438\u043c\u043d\u0438\u043d\u0433 \u043a\u043e\u043b\u0433\u0430\u043d 80 (\u0441\u0430\u043a\u0441\u043e\u043d) \u0444\u043e\u0438\u0437\u0438\u043d\u0438 \u0442\u0443\u043b\u0430\u0448\u0434\u0430 \u0418\u043d\u0432\u0435\u0441\u0442\u043e\u0440 \u0442\u0423\u043b\u043e\u0432\u043b\u0430\u0440\u043d\u0438 \u0438\u043a\u043a\u0438 \u0445\u0438\u043b \u0448\u0430\u043a\u043b\u0434\u0430 \u0430\u043c\u0430\u043b\u0433\u0430 \u043e\u0448\u0438\u0440\u0438\u0448\u0438 \u043c\u0443\u043c\u043a\u0438\u043d: \u2014 \u0443\u0437 \u043c\u0430\u0431\u043b\u0430\u0433\u043b\u0430\u0440\u0438 \u0445\u0438\u0441\u043e\u0431\u0438\u0434\u0430\u043d \u0451\u043a\u0438 \u2014 \u0416\u0430\u043c\u0433\u0430\u0440\u043c\u0430\u043d\u0438\u043d\u0433 \u0432\u0430\u043a\u043e\u043b\u0430\u0442\u043b\u0438 \u0442\u0438\u0436\u043e\u0440\u0430\u0442 \u0431\u0430\u043d\u043a\u0438 (\u043a\u0435\u0439\u0438\u043d\u0433\u0438 \u0423\u0440\u0438\u043d\u043b\u0430\u0440\u0434\u0430 \u2014 \u0411\u0430\u043d\u043a) \u043e\u0440\u043a\u0430\u043b\u0438 \u043a\u044a\u0448\u0438\u043c\u0447\u0430 \u043a\u0430\u0444\u0438\u043b\u043b\u0438\u043a \u0445\u0430\u0442\u0438 \u043e\u043b\u0438\u0448 \u0439\u0423\u043b\u0438 \u0431\u0438\u043b\u0430\u043d. \u041a\u044a\u0448\u0438\u043c\u0447\u0430 \u043a\u0430\u0444\u0438\u043b\u043b\u0438\u043a \u0445\u0430\u0442\u0438 \u0430\u0441\u043e\u0441\u0438\u0434\u0430 \u0416\u0430\u043c\u0433\u0430\u0440\u043c\u0430, \u0418\u043d\u0432\u0435\u0441\u0442\u043e\u0440 \u0432\u0430 \u0411\u0430\u043d\u043a \u0443\u0448\u0431\u0443 \u0448\u0430\u0440\u0442\u043d\u043e\u043c\u0430\u0434\u0430 \u043d\u0430\u0437\u0430\u0440\u0434\u0430 \u0442\u0443\u0442\u0438\u043b\u0433\u0430\u043d \u0445\u0430\u0440\u0430\u043a\u0430\u0442\u043b\u0430\u0440\u043d\u0438 \u0430\u043c\u0430\u043b\u0433\u0430 \u043e\u0448\u0438\u0440\u0430\u0434\u0438\u043b\u0430\u0440.\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM-interop script: adds \"Client Name\" / \"Company Name\" paragraphs\n# right after the document title and rewrites the body text of the \"2.2.\"\n# and \"2.3.\" clauses.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphIndexByPrefix {\n    param([string]$Prefix)\n    $i = 1\n    foreach ($p in $d.Paragraphs) {\n        $t = $p.Range.Text\n        $t = $t.TrimEnd([char]13, [char]7)\n        if ($t.StartsWith($Prefix)) {\n            return $i\n        }\n        $i++\n    }\n    return -1\n}\n\n# --- 1) Insert the two new paragraphs just before the \"1. \u0428\u0430\u0440\u0442\u043d\u043e\u043c\u0430 \u043f\u0440\u0435\u0434\u043c\u0435\u0442\u0438\"\n#        heading (i.e. right after the document title paragraph). ---\n$headingIdx = Get-ParagraphIndexByPrefix('1. \u0428\u0430\u0440\u0442\u043d\u043e\u043c\u0430 \u043f\u0440\u0435\u0434\u043c\u0435\u0442\u0438')\nif ($headingIdx -lt 0) {\n    throw \"Could not find the '1. \u0428\u0430\u0440\u0442\u043d\u043e\u043c\u0430 \u043f\u0440\u0435\u0434\u043c\u0435\u0442\u0438' heading paragraph\"\n}\n\n$headingPara = $d.Paragraphs.Item($headingIdx)\n$headingPara.Range.InsertParagraphBefore()\n$clientPara = $d.Paragraphs.Item($headingIdx)\n$clientPara.Range.Text = 'Client Name: Whitney Hurley'\n$clientPara.Style = $d.Styles.Item(\"Normal\")\n\n$headingPara = $d.Paragraphs.Item($headingIdx + 1)\n$headingPara.Range.InsertParagraphBefore()\n$companyPara = $d.Paragraphs.Item($headingIdx + 1)\n$companyPara.Range.Text = 'Company Name: Roberts and Richards Co'\n$companyPara.Style = $d.Styles.Item(\"Normal\")\n\n# --- 2) Rewrite the \"2.2.\" and \"2.3.\" paragraph bodies in place. ---\n$p22Idx = Get-ParagraphIndexByPrefix('2.2. \u0418\u043d\u0432\u0435\u0441\u0442\u043e\u0440 \u0439\u0438\u0433\u0438\u043c \u043c\u0438\u043a\u0434\u043e\u0440\u0438\u043d\u0438 \u0442\u0443\u043b\u0438\u043a \u0445\u0430\u0436\u043c\u0434\u0430')\nif ($p22Idx -lt 0) {\n    throw \"Could not find the '2.2.' paragraph\"\n}\n$p22 = $d.Paragraphs.Item($p22Idx)\n$p22.Range.Text = '2.2. \u0418\u043d\u0432\u0435\u0441\u0442\u043e\u0440 \u0439\u0438\u0433\u0438\u043c \u043c\u0438\u043a\u0434\u043e\u0440\u0438\u043d\u0438 \u0442\u0443\u043b\u0438\u043a \u0445\u0430\u0436\u043c\u0434\u0430 \u0416\u0430\u043c\u0433\u0430\u0440\u043c\u0430 \u0445\u0438\u0441\u043e\u0431 \u0440\u0430\u043a\u0430\u043c\u0438\u0433\u0430 \u0443\u0442\u043a\u0430\u0437\u0438\u0448 \u0439\u0423\u043b\u0438 \u0431\u0438\u043b\u0430\u043d \u043a\u0443\u0439\u0438\u0434\u0430\u0433\u0438 \u0442\u0430\u0440\u0442\u0438\u0431\u0434\u0430 \u0442\u0423\u043b\u0430\u043d\u0430\u0434\u0438: \u2014 20 (\u0439\u0438\u0433\u0438\u0440\u043c\u0430) \u0444\u043e\u0438\u0437\u0438 432.480.000 (\u0442\u0443\u0440\u0442 \u044e\u0437 \u0443\u0442\u0442\u0438\u0437 \u0438\u043a\u043a\u0438 \u043c\u0438\u043b\u043b\u0438\u043e\u043d \u0442\u0443\u0440\u0442 \u044e\u0437 \u0441\u0430\u043a\u0441\u043e\u043d \u043c\u0438\u043d\u0433) \u0441\u0443\u043c \u043c\u0438\u043a\u0434\u043e\u0440\u0438\u0434\u0430\u0433\u0438 \u043e\u043b\u0434\u0438\u043d\u0434\u0430\u043d \u0442\u0443\u043b\u043e\u0432\u043d\u0438 \u043c\u0430\u0437\u043a\u0443\u0440 \u0448\u0430\u0440\u0442\u043d\u043e\u043c\u0430 \u0438\u043c\u0437\u043e\u043b\u0430\u043d\u0433\u0430\u043d \u0441\u0430\u043d\u0430\u0434\u0430\u043d \u0431\u043e\u0448\u043b\u0430\u0431 \u0417 (\u0443\u0447) \u0438\u0448 \u043a\u0443\u043d\u0438\u0434\u0430\u043d \u043a\u0435\u0447\u0438\u043a\u0442\u0438\u0440\u043c\u0430\u0433\u0430\u043d \u0445\u043e\u043b\u0434\u0430 \u0430\u043c\u0430\u043b\u0433\u0430 \u043e\u0448\u0438\u0440\u0430\u0434\u0438; \u2014 \u0439\u0438\u0433\u0438\u043c\u043d\u0438\u043d\u0433 \u043a\u043e\u043b\u0433\u0430\u043d 80 (\u0441\u0430\u043a\u0441\u043e\u043d) \u0444\u043e\u0438\u0437\u0438 1.729.920.000 (\u0431\u0438\u0440 \u043c\u0438\u043b\u043b\u0438\u0430\u0440\u0434 \u0439\u0435\u0442\u0442\u0438 \u044e\u0437 \u0439\u0438\u0433\u0438\u0440\u043c\u0430 \u0442\u0443\u043a\u043a\u0438\u0437 \u043c\u0438\u043b\u043b\u0438\u043e\u043d \u0442\u0443\u043a\u043a\u0438\u0437 \u044e\u0437 \u0439\u0438\u0433\u0438\u0440\u043c\u0430 \u043c\u0438\u043d\u0433) \u0441\u0443\u043c \u043c\u0438\u043a\u0434\u043e\u0440\u0434\u0430\u0433\u0438 \u043a\u044a\u0448\u0438\u043c\u0447\u0430 \u0442\u0443\u043b\u043e\u0432\u043d\u0438 \u043a\u0443\u043f \u043a\u0430\u0432\u0430\u0442\u043b\u0438 \u0442\u0443\u0440\u0430\u0440 \u0436\u043e\u0439 \u0431\u0438\u043d\u043e\u0441\u0438 (\u0438\u043d\u0448\u043e\u043e\u0442\u0438) \u043a\u0443\u0440\u0438\u0431 \u0431\u0438\u0442\u043a\u0430\u0437\u0438\u043b\u0433\u0443\u043d\u0433\u0430 \u043a\u0430\u0434\u0430\u0440, \u0430\u043c\u043c\u043e 2023 \u0439\u0438\u043b 1 \u043d\u043e\u044f\u0431\u0440\u0434\u0430\u043d \u043a\u0435\u0447\u0438\u043a\u0442\u0438\u0440\u043c\u0430\u0441\u0434\u0430\u043d \u043a\u0443\u0439\u0438\u0434\u0430\u0433\u0438 \u043c\u0443\u0445\u043b\u0430\u0442\u043b\u0430\u0440\u0434\u0430 \u0430\u043c\u0430\u043b\u0433\u0430 \u043e\u0448\u0438\u0440\u0430\u0434\u0438.'\n\n$p23Idx = Get-ParagraphIndexByPrefix('2.3. \u041c\u0430\u0437\u043a\u0443\u0440 \u0448\u0430\u0440\u0442\u043d\u043e\u043c\u0430\u0434\u0430 \u043a\u0423\u0440\u0441\u0430\u0442\u0438\u043b\u0433\u0430\u043d')\nif ($p23Idx -lt 0) {\n    throw \"Could not find the '2.3.' paragraph\"\n}\n$p23 = $d.Paragraphs.Item($p23Idx)\n$p23.Range.Text = '2.3. \u0418\u0438\u0433\u0438\u043c\u043d\u0438\u043d\u0433 \u043a\u043e\u043b\u0433\u0430\u043d 80 (\u0441\u0430\u043a\u0441\u043e\u043d) \u0444\u043e\u0438\u0437\u0438\u043d\u0438 \u0442\u0443\u043b\u0430\u0448\u0434\u0430 \u0418\u043d\u0432\u0435\u0441\u0442\u043e\u0440 \u0442\u0423\u043b\u043e\u0432\u043b\u0430\u0440\u043d\u0438 \u0438\u043a\u043a\u0438 \u0445\u0438\u043b \u0448\u0430\u043a\u043b\u0434\u0430 \u0430\u043c\u0430\u043b\u0433\u0430 \u043e\u0448\u0438\u0440\u0438\u0448\u0438 \u043c\u0443\u043c\u043a\u0438\u043d: \u2014 \u0443\u0437 \u043c\u0430\u0431\u043b\u0430\u0433\u043b\u0430\u0440\u0438 \u0445\u0438\u0441\u043e\u0431\u0438\u0434\u0430\u043d \u0451\u043a\u0438 \u2014 \u0416\u0430\u043c\u0433\u0430\u0440\u043c\u0430\u043d\u0438\u043d\u0433 \u0432\u0430\u043a\u043e\u043b\u0430\u0442\u043b\u0438 \u0442\u0438\u0436\u043e\u0440\u0430\u0442 \u0431\u0430\u043d\u043a\u0438 (\u043a\u0435\u0439\u0438\u043d\u0433\u0438 \u0423\u0440\u0438\u043d\u043b\u0430\u0440\u0434\u0430 \u2014 \u0411\u0430\u043d\u043a) \u043e\u0440\u043a\u0430\u043b\u0438 \u043a\u044a\u0448\u0438\u043c\u0447\u0430 \u043a\u0430\u0444\u0438\u043b\u043b\u0438\u043a \u0445\u0430\u0442\u0438 \u043e\u043b\u0438\u0448 \u0439\u0423\u043b\u0438 \u0431\u0438\u043b\u0430\u043d. \u041a\u044a\u0448\u0438\u043c\u0447\u0430 \u043a\u0430\u0444\u0438\u043b\u043b\u0438\u043a \u0445\u0430\u0442\u0438 \u0430\u0441\u043e\u0441\u0438\u0434\u0430 \u0416\u0430\u043c\u0433\u0430\u0440\u043c\u0430, \u0418\u043d\u0432\u0435\u0441\u0442\u043e\u0440 \u0432\u0430 \u0411\u0430\u043d\u043a \u0443\u0448\u0431\u0443 \u0448\u0430\u0440\u0442\u043d\u043e\u043c\u0430\u0434\u0430 \u043d\u0430\u0437\u0430\u0440\u0434\u0430 \u0442\u0443\u0442\u0438\u043b\u0433\u0430\u043d \u0445\u0430\u0440\u0430\u043a\u0430\u0442\u043b\u0430\u0440\u043d\u0438 \u0430\u043c\u0430\u043b\u0433\u0430 \u043e\u0448\u0438\u0440\u0430\u0434\u0438\u043b\u0430\u0440.'\n"}
